# Rename the header row suffixes:
#   *_old -> *_FV2310   (columns A:J)
#   *_new -> *_FV2404   (columns L:U)
# and turn the data range A1:U60 into a real Excel Table ("Table1")
# with a frozen header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 21
$lastRow = 60

for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value()
    if ($val -like "*_old") {
        $cell.Value = ($val -replace "_old$", "_FV2310")
    } elseif ($val -like "*_new") {
        $cell.Value = ($val -replace "_new$", "_FV2404")
    }
}

# Turn the range into a native Excel table (ListObject) so the
# workbook gains xl/tables/table1.xml + the <tableParts> wiring.
$dataRange = $ws.Range("A1:U60")
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1, $null)
$tbl.Name = "Table1"

# Freeze the header row (split after row 1) and leave the selection
# parked on the first cell of the scrolling pane.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()
